# 2017-01-31 update: energy.gov - chunk 7
# Rolls the report forward from "October 2016"/"October 2015" YTD figures
# to "November 2016"/"November 2015" YTD figures (title, column headers,
# and the underlying state-level receipt data).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Table_4_08_B")

# --- Title (row 1) ---------------------------------------------------
$ws.Range("A1").Value = "Table 4.8.B. Receipts of Petroleum Coke Delivered for Electricity Generation by State, (Year-to-Date) November 2016 and 2015"

# --- Column headers (row 5): "October 20xx YTD" -> "November 20xx YTD"
foreach ($col in @("B", "E", "G", "I", "K")) {
    $ws.Range($col + "5").Value = "November 2016 YTD"
}
foreach ($col in @("C", "F", "H", "J", "L")) {
    $ws.Range($col + "5").Value = "November 2015 YTD"
}

# --- Data updates ------------------------------------------------------
# row 13 - Middle Atlantic
$ws.Range("C13").Value = 99
$ws.Range("L13").Value = 99

# row 16 - Pennsylvania
$ws.Range("C16").Value = 99
$ws.Range("L16").Value = 99

# row 17 - East North Central
$ws.Range("B17").Value = 940
$ws.Range("C17").Value = 1192
$ws.Range("D17").Value = -0.21
$ws.Range("E17").Value = 485
$ws.Range("F17").Value = 666
$ws.Range("G17").Value = 437
$ws.Range("H17").Value = 458
$ws.Range("L17").Value = 68

# row 19 - Indiana
$ws.Range("C19").Value = 353
$ws.Range("D19").Value = -0.54
$ws.Range("F19").Value = 353

# row 20 - Michigan
$ws.Range("B20").Value = 278
$ws.Range("C20").Value = 295
$ws.Range("D20").Value = -0.057000000000000002
$ws.Range("E20").Value = 278
$ws.Range("F20").Value = 282

# row 21 - Ohio
$ws.Range("B21").Value = 437
$ws.Range("C21").Value = 445
$ws.Range("D21").Value = -0.017999999999999999
$ws.Range("G21").Value = 437
$ws.Range("H21").Value = 445

# row 22 - Wisconsin
$ws.Range("B22").Value = 63
$ws.Range("C22").Value = 98
$ws.Range("D22").Value = -0.35
$ws.Range("E22").Value = 45
$ws.Range("F22").Value = 30
$ws.Range("L22").Value = 68

# row 31 - South Atlantic
$ws.Range("B31").Value = 1334
$ws.Range("C31").Value = 1043
$ws.Range("D31").Value = 0.28000000000000003
$ws.Range("E31").Value = 1241
$ws.Range("F31").Value = 936
$ws.Range("K31").Value = 93

# row 34 - Florida
$ws.Range("B34").Value = 1241
$ws.Range("C34").Value = 936
$ws.Range("D34").Value = 0.33
$ws.Range("E34").Value = 1241
$ws.Range("F34").Value = 936

# row 35 - Georgia
$ws.Range("B35").Value = 93
$ws.Range("D35").Value = -0.13
$ws.Range("K35").Value = 93

# row 41 - East South Central
$ws.Range("B41").Value = 77
$ws.Range("C41").Value = 559
$ws.Range("D41").Value = -0.86
$ws.Range("E41").Value = 77
$ws.Range("F41").Value = 559

# row 43 - Kentucky
$ws.Range("B43").Value = 77
$ws.Range("C43").Value = 559
$ws.Range("D43").Value = -0.86
$ws.Range("E43").Value = 77
$ws.Range("F43").Value = 559

# row 46 - West South Central
$ws.Range("B46").Value = 1451
$ws.Range("C46").Value = 1611
$ws.Range("D46").Value = -0.099000000000000005
$ws.Range("E46").Value = 1451
$ws.Range("F46").Value = 1611

# row 48 - Louisiana
$ws.Range("B48").Value = 1451
$ws.Range("C48").Value = 1611
$ws.Range("D48").Value = -0.099000000000000005
$ws.Range("E48").Value = 1451
$ws.Range("F48").Value = 1611

# row 67 - U.S. Total
$ws.Range("B67").Value = 3803
$ws.Range("C67").Value = 4504
$ws.Range("D67").Value = -0.16
$ws.Range("E67").Value = 3254
$ws.Range("F67").Value = 3772
$ws.Range("G67").Value = 437
$ws.Range("H67").Value = 458
$ws.Range("K67").Value = 112
$ws.Range("L67").Value = 274
